$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (paragraph 2, right after the title).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# 2. Insert a new bold "Play Almighty Sparta Free..." paragraph right before the
#    trailing "Prompt for DALLE" paragraph (now the last paragraph in the body).
$n = $d.Paragraphs.Count
$beforeLast = $d.Paragraphs.Item($n - 1)
$insertPos = $beforeLast.Range.End - 1
$insertPoint = $d.Range($insertPos, $insertPos)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Almighty Sparta Free: Game Review and Pros &amp; Cons</w:t></w:r></w:p>'
$insertPoint.InsertXML($newParaXml) | Out-Null

# 3. Replace the text of the final "Prompt for DALLE" paragraph with the meta
#    description text, keeping its existing (italic) run formatting.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $last.Range
$lastRange.MoveEnd(1, -1) | Out-Null
$lastRange.Text = "Read our review of Almighty Sparta online slot game. Learn about its pros & cons before you play for free. Discover if the game is worth your time."
